$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 11.39997876197891
$ws.Range("D2").Value = 4.86064255170043
$ws.Range("E2").Value = 11.2074708779952
$ws.Range("F2").Value = 57.7746890519082
$ws.Range("G2").Value = 78.59567557695948
$ws.Range("H2").Value = 26.17387131447306
$ws.Range("J2").Value = 10.57022822086932
$ws.Range("L2").Value = 8.200946204045596
$ws.Range("M2").Value = 35.29548741705593
$ws.Range("N2").Value = 17.51568045241979

# Row 3
$ws.Range("C3").Value = 11.44342066004473
$ws.Range("D3").Value = 4.769923155248835
$ws.Range("E3").Value = 11.23714557018159
$ws.Range("F3").Value = 57.60087961874182
$ws.Range("G3").Value = 78.05104785854728
$ws.Range("H3").Value = 26.15390854441009
$ws.Range("J3").Value = 10.61659357856541
$ws.Range("L3").Value = 8.18974324514293
$ws.Range("M3").Value = 34.66295100685373
$ws.Range("N3").Value = 17.34036893700447

# Row 4
$ws.Range("C4").Value = 11.47250053978795
$ws.Range("D4").Value = 4.712631826147979
$ws.Range("E4").Value = 11.25630219689479
$ws.Range("F4").Value = 57.51319000802609
$ws.Range("G4").Value = 77.74437839097132
$ws.Range("H4").Value = 26.14928184959506
$ws.Range("J4").Value = 10.64639258030702
$ws.Range("L4").Value = 8.182957995688476
$ws.Range("M4").Value = 34.27342482090467
$ws.Range("N4").Value = 17.2342335993989

# Row 5
$ws.Range("C5").Value = 11.48495367632318
$ws.Range("D5").Value = 4.688897017897072
$ws.Range("E5").Value = 11.26434489990705
$ws.Range("F5").Value = 57.48223816858652
$ws.Range("G5").Value = 77.62645438356881
$ws.Range("H5").Value = 26.149307367341
$ws.Range("J5").Value = 10.65887171208725
$ws.Range("L5").Value = 8.180216628149397
$ws.Range("M5").Value = 34.11459003613805
$ws.Range("N5").Value = 17.19140454046106

# Row 6
$ws.Range("C6").Value = 11.487057855494
$ws.Range("D6").Value = 4.684932776030601
$ws.Range("E6").Value = 11.26569467517481
$ws.Range("F6").Value = 57.47738727774685
$ws.Range("G6").Value = 77.60730059319953
$ws.Range("H6").Value = 26.14942676291883
$ws.Range("J6").Value = 10.66096417997059
$ws.Range("L6").Value = 8.179762861474764
$ws.Range("M6").Value = 34.0882151928612
$ws.Range("N6").Value = 17.18431954156186

# Row 7
$ws.Range("C7").Value = 11.47266604916223
$ws.Range("D7").Value = 4.712313286146421
$ws.Range("E7").Value = 11.25640970608278
$ws.Range("F7").Value = 57.51275322526883
$ws.Range("G7").Value = 77.7427594116474
$ws.Range("H7").Value = 26.14927446780713
$ws.Range("J7").Value = 10.64655951696297
$ws.Range("L7").Value = 8.182920928540293
$ws.Range("M7").Value = 34.27128287159784
$ws.Range("N7").Value = 17.23365422670986

# Row 8
$ws.Range("C8").Value = 11.41445672063245
$ws.Range("D8").Value = 4.829699075400315
$ws.Range("E8").Value = 11.21750887513623
$ws.Range("F8").Value = 57.71080998117041
$ws.Range("G8").Value = 78.40217955060686
$ws.Range("H8").Value = 26.16540016099187
$ws.Range("J8").Value = 10.58593960120637
$ws.Range("L8").Value = 8.197064147939084
$ws.Range("M8").Value = 35.07773876076436
$ws.Range("N8").Value = 17.45494633738203

# Row 9
$ws.Range("C9").Value = 11.31950443466837
$ws.Range("D9").Value = 5.046808912579527
$ws.Range("E9").Value = 11.14861655779205
$ws.Range("F9").Value = 58.25004535565651
$ws.Range("G9").Value = 79.91149431340851
$ws.Range("H9").Value = 26.25785369715861
$ws.Range("J9").Value = 10.47756389354351
$ws.Range("L9").Value = 8.22553290942818
$ws.Range("M9").Value = 36.64202859323315
$ws.Range("N9").Value = 17.89894323396332

# Row 10
$ws.Range("C10").Value = 11.26158841008832
$ws.Range("D10").Value = 5.197756875411887
$ws.Range("E10").Value = 11.10245612689682
$ws.Range("F10").Value = 58.73758999534569
$ws.Range("G10").Value = 81.14609457579155
$ws.Range("H10").Value = 26.36316261000065
$ws.Range("J10").Value = 10.40426041255933
$ws.Range("L10").Value = 8.24688815914738
$ws.Range("M10").Value = 37.77042243961448
$ws.Range("N10").Value = 18.22870289590215

# Row 11
$ws.Range("C11").Value = 11.23784456571821
$ws.Range("D11").Value = 5.264456272105621
$ws.Range("E11").Value = 11.082412752797
$ws.Range("F11").Value = 58.97897284796899
$ws.Range("G11").Value = 81.73349374357605
$ws.Range("H11").Value = 26.41921941788623
$ws.Range("J11").Value = 10.37226779066426
$ws.Range("L11").Value = 8.256696797461526
$ws.Range("M11").Value = 38.27714891856509
$ws.Range("N11").Value = 18.3789487527411

# Row 12
$ws.Range("C12").Value = 11.22923029795166
$ws.Range("D12").Value = 5.289422076694332
$ws.Range("E12").Value = 11.07495935983099
$ws.Range("F12").Value = 59.07316505824053
$ws.Range("G12").Value = 81.95948688135036
$ws.Range("H12").Value = 26.44161859917279
$ws.Range("J12").Value = 10.36034633533588
$ws.Range("L12").Value = 8.260424381034543
$ws.Range("M12").Value = 38.46792370261712
$ws.Range("N12").Value = 18.43583215150177

# Row 13
$ws.Range("C13").Value = 11.2310687222614
$ws.Range("D13").Value = 5.284058356643214
$ws.Range("E13").Value = 11.07655851829836
$ws.Range("F13").Value = 59.05275575471042
$ws.Range("G13").Value = 81.91065943182579
$ws.Range("H13").Value = 26.43674244788264
$ws.Range("J13").Value = 10.36290524867712
$ws.Range("L13").Value = 8.259620993401223
$ws.Range("M13").Value = 38.42688890094932
$ws.Range("N13").Value = 18.42358251808096

# Row 14
$ws.Range("C14").Value = 11.23712829279087
$ws.Range("D14").Value = 5.266516120658947
$ws.Range("E14").Value = 11.08179682442746
$ws.Range("F14").Value = 58.98666647129191
$ws.Range("G14").Value = 81.75201591079299
$ws.Range("H14").Value = 26.42103873804511
$ws.Range("J14").Value = 10.37128313493025
$ws.Range("L14").Value = 8.257003203644562
$ws.Range("M14").Value = 38.29286726274691
$ws.Range("N14").Value = 18.38362907112626

# Row 15
$ws.Range("C15").Value = 11.24088913588431
$ws.Range("D15").Value = 5.255732740250674
$ws.Range("E15").Value = 11.08502320456694
$ws.Range("F15").Value = 58.94654661421329
$ws.Range("G15").Value = 81.65530095854689
$ws.Range("H15").Value = 26.41157230928298
$ws.Range("J15").Value = 10.3764399910583
$ws.Range("L15").Value = 8.255401451620086
$ws.Range("M15").Value = 38.21062565202958
$ws.Range("N15").Value = 18.35915360497389

# Row 16
$ws.Range("C16").Value = 11.26319273095343
$ws.Range("D16").Value = 5.193357546711792
$ws.Range("E16").Value = 11.10378516502272
$ws.Range("F16").Value = 58.72220676694309
$ws.Range("G16").Value = 81.10821280136751
$ws.Range("H16").Value = 26.35966318756552
$ws.Range("J16").Value = 10.40637834093988
$ws.Range("L16").Value = 8.246249011936314
$ws.Range("M16").Value = 37.73715988871641
$ws.Range("N16").Value = 18.21888483188229

# Row 17
$ws.Range("C17").Value = 11.2775437799364
$ws.Range("D17").Value = 5.154581919214704
$ws.Range("E17").Value = 11.11553913547962
$ws.Range("F17").Value = 58.58957956133561
$ws.Range("G17").Value = 80.77909037068368
$ws.Range("H17").Value = 26.32990654085972
$ws.Range("J17").Value = 10.42509037490074
$ws.Range("L17").Value = 8.240658189717333
$ws.Range("M17").Value = 37.44489742381492
$ws.Range("N17").Value = 18.13286248312639

# Row 18
$ws.Range("C18").Value = 11.28604290947075
$ws.Range("D18").Value = 5.132094494124898
$ws.Range("E18").Value = 11.12238966557966
$ws.Range("F18").Value = 58.51514270000685
$ws.Range("G18").Value = 80.59222324505903
$ws.Range("H18").Value = 26.3135586699494
$ws.Range("J18").Value = 10.43598050973144
$ws.Range("L18").Value = 8.237451362828724
$ws.Range("M18").Value = 37.27618522286053
$ws.Range("N18").Value = 18.08340964445018

# Row 19
$ws.Range("C19").Value = 11.28896252671976
$ws.Range("D19").Value = 5.12444919540277
$ws.Range("E19").Value = 11.12472461085743
$ws.Range("F19").Value = 58.4902577335403
$ws.Range("G19").Value = 80.52937590342208
$ws.Range("H19").Value = 26.30815534476821
$ws.Range("J19").Value = 10.43968965082616
$ws.Range("L19").Value = 8.23636711679883
$ws.Range("M19").Value = 37.21896246276008
$ws.Range("N19").Value = 18.06667144416553

# Row 20
$ws.Range("C20").Value = 11.27599073185885
$ws.Range("D20").Value = 5.158728834047294
$ws.Range("E20").Value = 11.11427860057273
$ws.Range("F20").Value = 58.60350703210248
$ws.Range("G20").Value = 80.81387500870498
$ws.Range("H20").Value = 26.33299476554636
$ws.Range("J20").Value = 10.42308526304005
$ws.Range("L20").Value = 8.241252426501926
$ws.Range("M20").Value = 37.47607367274522
$ws.Range("N20").Value = 18.14201745946287

# Row 21
$ws.Range("C21").Value = 11.23533819570551
$ws.Range("D21").Value = 5.271676687698619
$ws.Range("E21").Value = 11.08025450580174
$ws.Range("F21").Value = 59.00600317531441
$ws.Range("G21").Value = 81.79851802438559
$ws.Range("H21").Value = 26.42561951065306
$ws.Range("J21").Value = 10.36881710565005
$ws.Range("L21").Value = 8.257771754017295
$ws.Range("M21").Value = 38.33226410082345
$ws.Range("N21").Value = 18.39536503593339

# Row 22
$ws.Range("C22").Value = 11.21096823948929
$ws.Range("D22").Value = 5.343790095950204
$ws.Range("E22").Value = 11.05881363489448
$ws.Range("F22").Value = 59.28527639836285
$ws.Range("G22").Value = 82.46270239572611
$ws.Range("H22").Value = 26.49298327282628
$ws.Range("J22").Value = 10.33447679613764
$ws.Range("L22").Value = 8.268645285050718
$ws.Range("M22").Value = 38.88528515544169
$ws.Range("N22").Value = 18.56085512312184

# Row 23
$ws.Range("C23").Value = 11.22377282153694
$ws.Range("D23").Value = 5.305460541198664
$ws.Range("E23").Value = 11.07018446408394
$ws.Range("F23").Value = 59.13475107602687
$ws.Range("G23").Value = 82.10637513168038
$ws.Range("H23").Value = 26.45640575360532
$ws.Range("J23").Value = 10.35270212599579
$ws.Range("L23").Value = 8.26283489991456
$ws.Range("M23").Value = 38.59077868988651
$ws.Range("N23").Value = 18.47255272967832

# Row 24
$ws.Range("C24").Value = 11.27669209101977
$ws.Range("D24").Value = 5.156854619930627
$ws.Range("E24").Value = 11.11484819862239
$ws.Range("F24").Value = 58.59720477532947
$ws.Range("G24").Value = 80.79814154377499
$ws.Range("H24").Value = 26.33159621309764
$ws.Range("J24").Value = 10.42399136182236
$ws.Range("L24").Value = 8.240983748947579
$ws.Range("M24").Value = 37.46198101668246
$ws.Range("N24").Value = 18.13787848324548

# Row 25
$ws.Range("C25").Value = 11.34312185314274
$ws.Range("D25").Value = 4.989538579549968
$ws.Range("E25").Value = 11.16646772301886
$ws.Range("F25").Value = 58.08805520604869
$ws.Range("G25").Value = 79.48061855062093
$ws.Range("H25").Value = 26.22629308591267
$ws.Range("J25").Value = 10.5057666880327
$ws.Range("L25").Value = 8.217753588981362
$ws.Range("M25").Value = 36.2217858231563
$ws.Range("N25").Value = 17.77802041980726
